$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.840.32"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.105.36"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'227.69"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").Value = "'62.33"
$ws.Range("E7").Value = "  +2.88%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +6.83%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'22.08"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "2.124.29"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "38.851.75"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'71.66"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "'228.00"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").Value = "'172.15"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'0.138"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").Value = "'19.36"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "'4.58"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("E34").Value = "  +12.21%  "
$ws.Range("D35").Value = "'4.74"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").Value = "'0.0619"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("D37").Value = "'2.40"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'18.13"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").Value = "'102.14"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "'0.0228"
$ws.Range("E42").Value = "  +3.93%  "
$ws.Range("D43").Value = "1.525.75"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  +8.33%  "
$ws.Range("D45").Value = "'2.80"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").Value = "'4.16"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "'2.97"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "2.303.05"
$ws.Range("E51").Value = "  +0.82%  "
